$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-03-20 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-21 Friday", 2) | Out-Null

# Update table cell values
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "78-78=0"
$t.Cell(1,2).Range.Text = "45-1=44"
$t.Cell(1,3).Range.Text = "4+38=42"
$t.Cell(1,4).Range.Text = "3+83=86"
$t.Cell(1,5).Range.Text = "83-42=41"
$t.Cell(2,1).Range.Text = "72-36=36"
$t.Cell(2,2).Range.Text = "54+30=84"
$t.Cell(2,3).Range.Text = "37-13=24"
$t.Cell(2,4).Range.Text = "60-29=31"
$t.Cell(2,5).Range.Text = "46-13=33"
$t.Cell(3,1).Range.Text = "82-1=81"
$t.Cell(3,2).Range.Text = "91-65=26"
$t.Cell(3,3).Range.Text = "57-22=35"
$t.Cell(3,4).Range.Text = "75-7=68"
$t.Cell(3,5).Range.Text = "69+10=79"
$t.Cell(4,1).Range.Text = "42+0=42"
$t.Cell(4,2).Range.Text = "16-12=4"
$t.Cell(4,3).Range.Text = "39+2=41"
$t.Cell(4,4).Range.Text = "78-49=29"
$t.Cell(4,5).Range.Text = "77-52=25"
$t.Cell(5,1).Range.Text = "72+10=82"
$t.Cell(5,2).Range.Text = "47-3=44"
$t.Cell(5,3).Range.Text = "9+6=15"
$t.Cell(5,4).Range.Text = "31+22=53"
$t.Cell(5,5).Range.Text = "63+6=69"
$t.Cell(6,1).Range.Text = "60+16=76"
$t.Cell(6,2).Range.Text = "77-25=52"
$t.Cell(6,3).Range.Text = "78-47=31"
$t.Cell(6,4).Range.Text = "30+28=58"
$t.Cell(6,5).Range.Text = "24+61=85"
$t.Cell(7,1).Range.Text = "32-20=12"
$t.Cell(7,2).Range.Text = "28-19=9"
$t.Cell(7,3).Range.Text = "89-64=25"
$t.Cell(7,4).Range.Text = "57-34=23"
$t.Cell(7,5).Range.Text = "14+72=86"
$t.Cell(8,1).Range.Text = "23+53=76"
$t.Cell(8,2).Range.Text = "17+16=33"
$t.Cell(8,3).Range.Text = "50+45=95"
$t.Cell(8,4).Range.Text = "30+46=76"
$t.Cell(8,5).Range.Text = "17+56=73"
$t.Cell(9,1).Range.Text = "42+15=57"
$t.Cell(9,2).Range.Text = "13+32=45"
$t.Cell(9,3).Range.Text = "74-37=37"
$t.Cell(9,4).Range.Text = "20+57=77"
$t.Cell(9,5).Range.Text = "55+39=94"
$t.Cell(10,1).Range.Text = "45+17=62"
$t.Cell(10,2).Range.Text = "55+19=74"
$t.Cell(10,3).Range.Text = "7-4=3"
$t.Cell(10,4).Range.Text = "51-13=38"
$t.Cell(10,5).Range.Text = "50+10=60"
$t.Cell(11,1).Range.Text = "42-39=3"
$t.Cell(11,2).Range.Text = "3+76=79"
$t.Cell(11,3).Range.Text = "50-18=32"
$t.Cell(11,4).Range.Text = "39+14=53"
$t.Cell(11,5).Range.Text = "31+63=94"
$t.Cell(12,1).Range.Text = "27+17=44"
$t.Cell(12,2).Range.Text = "39-31=8"
$t.Cell(12,3).Range.Text = "65-1=64"
$t.Cell(12,4).Range.Text = "59-3=56"
$t.Cell(12,5).Range.Text = "95-42=53"
$t.Cell(13,1).Range.Text = "42+5=47"
$t.Cell(13,2).Range.Text = "33+34=67"
$t.Cell(13,3).Range.Text = "81-49=32"
$t.Cell(13,4).Range.Text = "54-46=8"
$t.Cell(13,5).Range.Text = "70-62=8"
$t.Cell(14,1).Range.Text = "49-47=2"
$t.Cell(14,2).Range.Text = "82-5=77"
$t.Cell(14,3).Range.Text = "34+41=75"
$t.Cell(14,4).Range.Text = "49-42=7"
$t.Cell(14,5).Range.Text = "11+56=67"
$t.Cell(15,1).Range.Text = "28+31=59"
$t.Cell(15,2).Range.Text = "60-56=4"
$t.Cell(15,3).Range.Text = "76-18=58"
$t.Cell(15,4).Range.Text = "87-65=22"
$t.Cell(15,5).Range.Text = "8+17=25"
$t.Cell(16,1).Range.Text = "30+15=45"
$t.Cell(16,2).Range.Text = "96-46=50"
$t.Cell(16,3).Range.Text = "97-57=40"
$t.Cell(16,4).Range.Text = "83-75=8"
$t.Cell(16,5).Range.Text = "87-23=64"
$t.Cell(17,1).Range.Text = "22+2=24"
$t.Cell(17,2).Range.Text = "62-10=52"
$t.Cell(17,3).Range.Text = "78+1=79"
$t.Cell(17,4).Range.Text = "76-45=31"
$t.Cell(17,5).Range.Text = "72-56=16"
$t.Cell(18,1).Range.Text = "46+33=79"
$t.Cell(18,2).Range.Text = "19+24=43"
$t.Cell(18,3).Range.Text = "35-3=32"
$t.Cell(18,4).Range.Text = "81-49=32"
$t.Cell(18,5).Range.Text = "90-70=20"
$t.Cell(19,1).Range.Text = "35+59=94"
$t.Cell(19,2).Range.Text = "40-9=31"
$t.Cell(19,3).Range.Text = "59-48=11"
$t.Cell(19,4).Range.Text = "48+34=82"
$t.Cell(19,5).Range.Text = "95-69=26"
$t.Cell(20,1).Range.Text = "66+11=77"
$t.Cell(20,2).Range.Text = "34+41=75"
$t.Cell(20,3).Range.Text = "99-76=23"
$t.Cell(20,4).Range.Text = "59+5=64"
$t.Cell(20,5).Range.Text = "33+26=59"
